# Updates cryptocurrency price (D) and 1h volume-change (E) columns
# to the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.655.54'
$ws.Range("E2").Value = '  +8.51%  '
$ws.Range("D3").Value = '2.508.11'
$ws.Range("E3").Value = '  +10.90%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.36%  '
$ws.Range("D5").Value = '''486.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +14.73%  '
$ws.Range("D6").Value = '''142.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +20.43%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '''0.513'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.64%  '
$ws.Range("D9").Value = '2.503.28'
$ws.Range("E9").Value = '  +10.27%  '
$ws.Range("D10").Value = '''0.0991'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.89%  '
$ws.Range("D11").Value = '''5.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.96%  '
$ws.Range("D12").Value = '''0.329'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.42%  '
$ws.Range("E13").Value = '  +2.06%  '
$ws.Range("D14").Value = '2.931.31'
$ws.Range("E14").Value = '  +10.20%  '
$ws.Range("D15").Value = '55.635.17'
$ws.Range("E15").Value = '  +8.27%  '
$ws.Range("D16").Value = '''20.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +11.12%  '
$ws.Range("D17").Value = '''0.0000138'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +19.49%  '
$ws.Range("D18").Value = '2.512.42'
$ws.Range("E18").Value = '  +9.82%  '
$ws.Range("D19").Value = '''4.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +13.40%  '
$ws.Range("D20").Value = '''320.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.43%  '
$ws.Range("D21").Value = '''10.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +14.00%  '
$ws.Range("D22").Value = '''0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '''5.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.66%  '
$ws.Range("D24").Value = '''58.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.48%  '
$ws.Range("D25").Value = '''0.169'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +15.84%  '
$ws.Range("D26").Value = '''0.410'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +14.71%  '
$ws.Range("D27").Value = '''1.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '2.626.35'
$ws.Range("E28").Value = '  +10.51%  '
$ws.Range("D29").Value = '''7.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.27%  '
$ws.Range("D30").Value = '0.0₃0802'
$ws.Range("E30").Value = '  +18.29%  '
$ws.Range("D31").Value = '''0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '''149.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.18%  '
$ws.Range("D33").Value = '''18.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.09%  '
$ws.Range("D34").Value = '''1.49'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.54%  '
$ws.Range("D35").Value = '''5.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.49%  '
$ws.Range("D36").Value = '''0.889'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.34%  '
$ws.Range("D37").Value = '''3.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.98%  '
$ws.Range("E38").Value = '  +15.61%  '
$ws.Range("D39").Value = '''34.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.06%  '
$ws.Range("D40").Value = '''0.614'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +20.69%  '
$ws.Range("D41").Value = '''0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").Value = '''0.0556'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.51%  '
$ws.Range("D43").Value = '''3.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.03%  '
$ws.Range("D44").Value = '''1.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.28%  '
$ws.Range("D45").Value = '2.004.92'
$ws.Range("E45").Value = '  +8.51%  '
$ws.Range("D46").Value = '''4.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +29.44%  '
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").Value = '''0.0909'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.67%  '
$ws.Range("E49").Value = '  +11.50%  '
$ws.Range("D50").Value = '''253.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +39.94%  '
$ws.Range("D51").Value = '''17.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.86%  '

Write-Output "Updated 96 cells in columns D and E (rows 2-51)"
